# Apply the edit described by the diff:
#  - Insert two new data rows into the "Hortaliza, Vega Modelo de Temuco -
#    Cebollín" weekly log, right after the existing row 607 (i.e. at rows
#    608-609), shifting the rest of the table down by two rows.
#  - Populate the two new rows with their new weekly observations.
# This naturally grows the sheet's used range from A1:R708 to A1:R710,
# matching the target dimension.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 608-609; everything from old row 608 onward
# shifts down to 610 onward (so old 707/708 become new 709/710).
$ws.Range("A608:A609").EntireRow.Insert()

# New row 608
$ws.Cells.Item(608, 1).Value = 10
$ws.Cells.Item(608, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(608, 3).Value = "La Araucanía"
$ws.Cells.Item(608, 4).Value = 45180
$ws.Cells.Item(608, 5).Value = 9
$ws.Cells.Item(608, 6).Value = 100112037
$ws.Cells.Item(608, 7).Value = "Cebollín"
$ws.Cells.Item(608, 8).Value = "Sin especificar"
$ws.Cells.Item(608, 9).Value = "Primera"
$ws.Cells.Item(608, 10).Value = 80
$ws.Cells.Item(608, 11).Value = 9000
$ws.Cells.Item(608, 12).Value = 9000
$ws.Cells.Item(608, 13).Value = 9000
$ws.Cells.Item(608, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(608, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(608, 16).Value = 750
$ws.Cells.Item(608, 17).Value = 12
$ws.Cells.Item(608, 18).Value = "Hortaliza"

# New row 609
$ws.Cells.Item(609, 1).Value = 10
$ws.Cells.Item(609, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(609, 3).Value = "La Araucanía"
$ws.Cells.Item(609, 4).Value = 45180
$ws.Cells.Item(609, 5).Value = 9
$ws.Cells.Item(609, 6).Value = 100112037
$ws.Cells.Item(609, 7).Value = "Cebollín"
$ws.Cells.Item(609, 8).Value = "Sin especificar"
$ws.Cells.Item(609, 9).Value = "Primera"
$ws.Cells.Item(609, 10).Value = 120
$ws.Cells.Item(609, 11).Value = 7000
$ws.Cells.Item(609, 12).Value = 7000
$ws.Cells.Item(609, 13).Value = 7000
$ws.Cells.Item(609, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(609, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(609, 16).Value = 583
$ws.Cells.Item(609, 17).Value = 12
$ws.Cells.Item(609, 18).Value = "Hortaliza"
